$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) from 45598 to 45599 for rows 2-32
for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45598) {
        $cell.Value2 = 45599
    }
}
